# #129 adding support of named ranges for DataModel class
#
# Update the ISEVEN formula in D12 to reference the existing named range
# "hjk" (Sheet1!$D$3, currently a blank cell) instead of the literal cell
# reference A1, and move the active selection to D1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D3 is blank, so ISEVEN(hjk) evaluates to TRUE (blank/0 is even).
$ws.Range("D12").Formula = "=ISEVEN(hjk)"

# Move the active cell/selection from D12 to D1.
$ws.Range("D1").Select()
